$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings (e.g. "1.011")
# are stored as text, matching the source inlineStr cells, not coerced to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "24.620.82"
$ws.Range("E2").Value = "  -0.73%  "

$ws.Range("D3").Value = "1.698.71"
$ws.Range("E3").Value = "  -0.36%  "

$ws.Range("D4").Value = "1.011"
$ws.Range("E4").Value = "  +1.40%  "

$ws.Range("D5").Value = "315.24"
$ws.Range("E5").Value = "  -0.28%  "

$ws.Range("E6").Value = "  +1.10%  "

$ws.Range("D7").Value = "0.3940"
$ws.Range("E7").Value = "  -0.40%  "

$ws.Range("D8").Value = "0.4067"
$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("D9").Value = "1.516"
$ws.Range("E9").Value = "  +1.24%  "

$ws.Range("D10").Value = "1.011"
$ws.Range("E10").Value = "  +1.45%  "

$ws.Range("D11").Value = "52.52"
$ws.Range("E11").Value = "  -0.60%  "

$ws.Range("D12").Value = "0.08796"
$ws.Range("E12").Value = "  -1.11%  "

$ws.Range("D13").Value = "7.591"
$ws.Range("E13").Value = "  +5.09%  "

$ws.Range("D14").Value = "24.73"
$ws.Range("E14").Value = "  +4.65%  "

$ws.Range("D15").Value = "0.00001364"
$ws.Range("E15").Value = "  +2.78%  "

$ws.Range("D16").Value = "8.022"
$ws.Range("E16").Value = "  -1.27%  "

$ws.Range("D17").Value = "1.693.24"
$ws.Range("E17").Value = "  -0.38%  "

$ws.Range("D18").Value = "99.06"
$ws.Range("E18").Value = "  -0.96%  "

$ws.Range("D19").Value = "0.07117"
$ws.Range("E19").Value = "  +1.57%  "

$ws.Range("D20").Value = "19.85"
$ws.Range("E20").Value = "  +0.54%  "

$ws.Range("D21").Value = "7.389"
$ws.Range("E21").Value = "  +4.76%  "

$ws.Range("D22").Value = "1.012"
$ws.Range("E22").Value = "  +1.26%  "

$ws.Range("D23").Value = "14.36"
$ws.Range("E23").Value = "  -0.43%  "

$ws.Range("D24").Value = "24.612.07"
$ws.Range("E24").Value = "  -0.67%  "

$ws.Range("D25").Value = "3.054"
$ws.Range("E25").Value = "  -6.68%  "

$ws.Range("D26").Value = "2.348"
$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").Value = "22.82"
$ws.Range("E27").Value = "  +0.09%  "

$ws.Range("D28").Value = "164.90"
$ws.Range("E28").Value = "  +1.24%  "

$ws.Range("D29").Value = "8.441"
$ws.Range("E29").Value = "  +12.83%  "

$ws.Range("D30").Value = "138.65"
$ws.Range("E30").Value = "  +1.80%  "

$ws.Range("D31").Value = "5.233"
$ws.Range("E31").Value = "  +0.85%  "

$ws.Range("D32").Value = "1.882.55"
$ws.Range("E32").Value = "  +0.12%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.08827"
$ws.Range("E33").Value = "  +2.37%  "

$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "7.511"
$ws.Range("E34").Value = "  +4.53%  "

$ws.Range("D35").Value = "1.045"
$ws.Range("E35").Value = "  -2.07%  "

$ws.Range("D36").Value = "1.994"
$ws.Range("E36").Value = "  +3.40%  "

$ws.Range("D37").Value = "0.2732"
$ws.Range("E37").Value = "  -1.03%  "

$ws.Range("D38").Value = "0.02889"
$ws.Range("E38").Value = "  +5.40%  "

$ws.Range("D39").Value = "10.86"
$ws.Range("E39").Value = "  -6.71%  "

$ws.Range("D40").Value = "14.35"
$ws.Range("E40").Value = "  -1.92%  "

$ws.Range("D41").Value = "0.09142"
$ws.Range("E41").Value = "  -0.73%  "

$ws.Range("D42").Value = "0.7840"
$ws.Range("E42").Value = "  +1.78%  "

$ws.Range("D43").Value = "1.468"
$ws.Range("E43").Value = "  -0.46%  "

$ws.Range("D44").Value = "16.66"
$ws.Range("E44").Value = "  +3.26%  "

$ws.Range("D45").Value = "0.7208"
$ws.Range("E45").Value = "  -0.17%  "

$ws.Range("D46").Value = "2.582"
$ws.Range("E46").Value = "  -0.46%  "

$ws.Range("D47").Value = "4.227"
$ws.Range("E47").Value = "  +0.15%  "

$ws.Range("D48").Value = "1.010"
$ws.Range("E48").Value = "  +0.86%  "

$ws.Range("D49").Value = "1.322"
$ws.Range("E49").Value = "  -1.01%  "

$ws.Range("D50").Value = "139.69"

$ws.Range("D51").Value = "91.99"
$ws.Range("E51").Value = "  +2.08%  "

# Reset column D back to the default (unstyled) cell format now that the
# values are committed as text, so no stray number-format style lingers.
$ws.Range("D2:D51").Style = "Normal"
